$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Populate shared strings in the same order they appear in the target file:
# 32 = "Шестерня планетарная", 33 = "planetary gear 2.jpg", 34 = "planetary gear 1.jpg"
$ws.Range("B14").Value = "Шестерня планетарная"
$ws.Range("C15").Value = "planetary gear 2.jpg"
$ws.Range("C14").Value = "planetary gear 1.jpg"

$ws.Range("A14").Value = "Реверсивный инжениринг"
$ws.Range("A15").Value = "Реверсивный инжениринг"
$ws.Range("B15").Value = "Шестерня планетарная"

# Row heights to match rows 8-13 pattern
$ws.Rows.Item(14).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 15.75

# Copy formats: A14/A15 and B14 get style "s=3" (same as A13/B13); B15/C14/C15 stay default
$ws.Range("A13").Copy()
$ws.Range("A14:A15").PasteSpecial(-4122)

$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("D12").Select()
